$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (avoids Excel auto-converting
# numeric-looking strings like "214.21" into real numbers), then reset
# the cell style back to Normal so no stray quote-prefix formatting sticks.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '27.085.34'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '1.623.48'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  -0.14%  '
Set-TextValue $ws.Range("D5") '214.21'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.49%  '
Set-TextValue $ws.Range("D9") '0.0626'
$ws.Range("E9").Value = '  -0.12%  '
Set-TextValue $ws.Range("D10") '20.29'
$ws.Range("E10").Value = '  +1.44%  '
Set-TextValue $ws.Range("D11") '0.0845'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '1.625.10'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("E13").Value = '  -0.54%  '
Set-TextValue $ws.Range("D14") '0.542'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").Value = '27.055.18'
$ws.Range("E15").Value = '  -0.55%  '
Set-TextValue $ws.Range("D16") '64.44'
$ws.Range("E16").Value = '  -4.40%  '
$ws.Range("E17").Value = '  +0.53%  '
Set-TextValue $ws.Range("D18") '215.80'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("E19").Value = '  -0.11%  '
Set-TextValue $ws.Range("D20") '6.90'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("E21").Value = '  -0.75%  '
Set-TextValue $ws.Range("D22") '2.42'
$ws.Range("E22").Value = '  -5.02%  '
Set-TextValue $ws.Range("D23") '9.04'
$ws.Range("E23").Value = '  -1.66%  '
Set-TextValue $ws.Range("D24") '147.07'
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("E27").Value = '  -0.54%  '
Set-TextValue $ws.Range("D28") '15.58'
$ws.Range("E28").Value = '  -1.12%  '
Set-TextValue $ws.Range("D29") '0.0505'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  -1.09%  '
Set-TextValue $ws.Range("D31") '3.35'
$ws.Range("E31").Value = '  -0.82%  '
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("D33").Value = '1.340.54'
$ws.Range("E33").Value = '  +6.29%  '
Set-TextValue $ws.Range("D34") '1.57'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("E36").Value = '  -0.80%  '
Set-TextValue $ws.Range("D37") '0.546'
$ws.Range("E37").Value = '  +0.30%  '
Set-TextValue $ws.Range("D38") '0.851'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  -0.15%  '
Set-TextValue $ws.Range("D40") '0.803'
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D41") '65.40'
$ws.Range("E41").Value = '  +5.83%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D42") '2.23'
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").Value = '1.760.14'
$ws.Range("E43").Value = '  -1.36%  '
Set-TextValue $ws.Range("D44") '5.22'
$ws.Range("E44").Value = '  -1.55%  '
Set-TextValue $ws.Range("D45") '90.40'
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D46") '1.61'
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D47") '0.853'
$ws.Range("E47").Value = '  +28.62%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  -0.34%  '
Set-TextValue $ws.Range("D50") '0.0994'
$ws.Range("E50").Value = '  +2.13%  '
Set-TextValue $ws.Range("D51") '7.57'
$ws.Range("E51").Value = '  -0.92%  '
